# Applies the "cleaned up completed tasks" edit:
#  - Removes the finished top task list, the "Post production" sub-notes
#    that were completed, and the now-irrelevant "Individual work load" /
#    "by percent?" block.
#  - Keeps "Post production notes, mid-stream changes" + the still-open
#    "Union and project implementations" item.
#  - Extends the "-From git" dev-log item with " include link on report"
#    and clears the stray "-Updates?" placeholder text (leaving its tab).
#  - Clears the finished "-Now is the time to complain or not." item
#    (leaving its tab) and relocates the _GoBack bookmark there, matching
#    where Word would leave it after the last edit.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Delete whole paragraphs that are no longer needed.  Working from the
#    highest paragraph index down keeps the not-yet-processed indices
#    stable as each Range.Delete() removes its paragraph mark too.
# ---------------------------------------------------------------------
$paragraphsToRemove = @(
    18,           # blank line right after the "-Updates?" placeholder
    14, 13, 12,   # "Individual work load distribution...", "-by percent? Or by Volume?", blank line
    9, 8,         # "-We removed Earning Report Table...", "-Design considerations for the DB App..."
    6, 5, 4, 3, 2, 1   # blank line, and the four completed top-of-doc tasks + its title
)
foreach ($idx in $paragraphsToRemove) {
    $d.Paragraphs.Item($idx).Range.Delete()
}

# ---------------------------------------------------------------------
# 2. "-From git" dev-log entry: append " include link on report" right
#    before the paragraph mark.
# ---------------------------------------------------------------------
$fromGit = $d.Paragraphs.Item(5).Range
$fromGit.End = $fromGit.End - 1
$fromGit.InsertAfter(" include link on report")

# ---------------------------------------------------------------------
# 3. "-Updates?" placeholder paragraph becomes just its leading tab.
#    Remove the "-Updates?" text between the two existing tab runs...
# ---------------------------------------------------------------------
$updates = $d.Paragraphs.Item(6).Range
$d.Range($updates.Start + 1, $updates.End - 2).Delete()
# ...then drop the now-bare leading tab character (it got folded into a
# plain text run by the delete above) so only the original trailing
# <w:tab/> run remains.
$updates = $d.Paragraphs.Item(6).Range
$d.Range($updates.Start, $updates.Start + 1).Delete()

# ---------------------------------------------------------------------
# 4. Final "-Now is the time to complain or not." item becomes just its
#    leading tab, and the _GoBack bookmark (originally sitting in the
#    "Individual work load..." paragraph we removed above) is re-created
#    here, mirroring where Word leaves it after the final edit.
# ---------------------------------------------------------------------
$complain = $d.Paragraphs.Item(14).Range
$d.Range($complain.Start + 1, $complain.End - 1).Delete()
$complain = $d.Paragraphs.Item(14).Range
$bookmarkRange = $d.Range($complain.End - 1, $complain.End - 1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
